# Product lookups now resolve an id/name for the article (instead of a bare
# "1") and the matching FACT/REMIS article description gets attached to each
# row on both the "Pendientes" sheet and carried through to "Facturados".
$wb = $excel.ActiveWorkbook
$wsPendientes = $wb.Worksheets.Item("Pendientes")
$wsFacturados = $wb.Worksheets.Item("Facturados")

# Row 2 - Luis Rendón / Fibra pp ... -> product "Pedrito" / "Clavos 100gr"
# Row 3 - Soga cabuya ... -> product "Juanito" / "Clavos 200gr"
$wsPendientes.Range("G2").Value = "Pedrito"
$wsPendientes.Range("G3").Value = "Juanito"
$wsPendientes.Range("O2").Value = "Clavos 100gr"
$wsPendientes.Range("O3").Value = "Clavos 200gr"

$wsPendientes.Range("R2").Value = 2000
$wsPendientes.Range("R3").Value = 3000

# Restore the view/selection on "Facturados" first so it no longer owns the
# active tab once we come back to "Pendientes" below.
$wsFacturados.Range("F5").Select() | Out-Null

# "Pendientes" becomes the active sheet again, scrolled over to column L with
# the selection parked on P5.
$wsPendientes.Activate() | Out-Null
$excel.ActiveWindow.ScrollColumn = 12
$wsPendientes.Range("P5").Select() | Out-Null
